$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before current row 460, shifting all data
# from row 460 downward down by two rows (rows 460-538 become 462-540).
$ws.Range("A460:A461").EntireRow.Insert()

# Populate the two newly-inserted rows with their data. Columns A, B, C,
# E, F, G, H, I, O, Q, R mirror the block they belong to (same market /
# category metadata as the surrounding "Ajo" - "Chino" - "Primera" rows),
# while D, J, K, L, M, N, P carry the new reported values.

# New row 460
$ws.Cells.Item(460, 1).Value2 = 8
$ws.Cells.Item(460, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(460, 3).Value2 = "Coquimbo"
$ws.Cells.Item(460, 4).Value2 = 45180
$ws.Cells.Item(460, 5).Value2 = 4
$ws.Cells.Item(460, 6).Value2 = 100112003
$ws.Cells.Item(460, 7).Value2 = "Ajo"
$ws.Cells.Item(460, 8).Value2 = "Chino"
$ws.Cells.Item(460, 9).Value2 = "Primera"
$ws.Cells.Item(460, 10).Value2 = 320
$ws.Cells.Item(460, 11).Value2 = 21500
$ws.Cells.Item(460, 12).Value2 = 22000
$ws.Cells.Item(460, 13).Value2 = 21750
$ws.Cells.Item(460, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(460, 15).Value2 = "China"
$ws.Cells.Item(460, 16).Value2 = 2175
$ws.Cells.Item(460, 17).Value2 = 10
$ws.Cells.Item(460, 18).Value2 = "Hortaliza"

# New row 461
$ws.Cells.Item(461, 1).Value2 = 8
$ws.Cells.Item(461, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(461, 3).Value2 = "Coquimbo"
$ws.Cells.Item(461, 4).Value2 = 45180
$ws.Cells.Item(461, 5).Value2 = 4
$ws.Cells.Item(461, 6).Value2 = 100112003
$ws.Cells.Item(461, 7).Value2 = "Ajo"
$ws.Cells.Item(461, 8).Value2 = "Chino"
$ws.Cells.Item(461, 9).Value2 = "Primera"
$ws.Cells.Item(461, 10).Value2 = 300
$ws.Cells.Item(461, 11).Value2 = 23500
$ws.Cells.Item(461, 12).Value2 = 24000
$ws.Cells.Item(461, 13).Value2 = 23750
$ws.Cells.Item(461, 14).Value2 = "$/malla 10 kilos"
$ws.Cells.Item(461, 15).Value2 = "China"
$ws.Cells.Item(461, 16).Value2 = 2375
$ws.Cells.Item(461, 17).Value2 = 10
$ws.Cells.Item(461, 18).Value2 = "Hortaliza"

# Make sure the date cells keep the same number format as the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(460, 4).NumberFormat = $ws.Cells.Item(462, 4).NumberFormat
$ws.Cells.Item(461, 4).NumberFormat = $ws.Cells.Item(462, 4).NumberFormat
